# The sheet contains a weekly price list for "Vega Modelo de Temuco - Brócoli".
# A new weekly record needs to be inserted at row 203 (pushing the existing
# rows 203-305 down to 204-306), growing the sheet from 305 to 306 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 203; Excel shifts rows
# 203:305 down to 204:306 and extends the used range accordingly.
$ws.Rows("203:203").Insert()

# Populate the newly inserted row 203 with the new weekly entry.
$ws.Range("A203").Value = 10
$ws.Range("B203").Value = 'Vega Modelo de Temuco'
$ws.Range("C203").Value = 'La Araucanía'
$ws.Range("D203").Value = 44523
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = 100112023
$ws.Range("G203").Value = 'Brócoli'
$ws.Range("H203").Value = 'Sin especificar'
$ws.Range("I203").Value = 'Primera'
$ws.Range("J203").Value = 3000
$ws.Range("K203").Value = 800
$ws.Range("L203").Value = 800
$ws.Range("M203").Value = 800
$ws.Range("N203").Value = '$/unidad'
$ws.Range("O203").Value = 'Región Metropolitana'
$ws.Range("P203").Value = 800
$ws.Range("Q203").Value = 1
$ws.Range("R203").Value = 'Hortaliza'
